$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3991.7856
$ws.Range("I62").Value = 5121.6665
$ws.Range("J62").Value = 1958
$ws.Range("K62").Value = 5121.6665
$ws.Range("L62").Value = 1958
$ws.Range("M62").Value = -4497.6665
$ws.Range("N62").Value = -3206

$ws.Range("H65").Value = 3991.7856
$ws.Range("I65").Value = 5121.6665
$ws.Range("J65").Value = 1958
$ws.Range("K65").Value = 25608.3325
$ws.Range("L65").Value = 9790
$ws.Range("M65").Value = -22488.3325
$ws.Range("N65").Value = -16030

$ws.Range("H70").Value = 1138.7826
$ws.Range("I70").Value = 1052.4667
$ws.Range("J70").Value = 1300.625
$ws.Range("K70").Value = 3157.4001
$ws.Range("L70").Value = 3901.875
$ws.Range("M70").Value = -2887.4001
$ws.Range("N70").Value = -4441.875

$ws.Range("H73").Value = 1138.7826
$ws.Range("I73").Value = 1052.4667
$ws.Range("J73").Value = 1300.625
$ws.Range("K73").Value = 3157.4001
$ws.Range("L73").Value = 3901.875
$ws.Range("M73").Value = -2221.4001
$ws.Range("N73").Value = -5773.875

$ws.Range("H96").Value = 83342824
$ws.Range("I96").Value = 4465.5
$ws.Range("K96").Value = 13396.5
$ws.Range("M96").Value = -12023.5

$ws.Range("H98").Value = 37008.953
$ws.Range("I98").Value = 1431.8462
$ws.Range("J98").Value = 88398.11
$ws.Range("K98").Value = 1431.8462
$ws.Range("L98").Value = 88398.11
$ws.Range("M98").Value = 66.15380000000005
$ws.Range("N98").Value = -91394.11

$ws.Range("H100").Value = 2522.9412
$ws.Range("I100").Value = 3260
$ws.Range("J100").Value = 2215.8333
$ws.Range("K100").Value = 3260
$ws.Range("L100").Value = 2215.8333
$ws.Range("M100").Value = -2719
$ws.Range("N100").Value = -3297.8333

$ws.Range("H122").Value = 37008.953
$ws.Range("I122").Value = 1431.8462
$ws.Range("J122").Value = 88398.11
$ws.Range("K122").Value = 4295.5386
$ws.Range("L122").Value = 265194.33
$ws.Range("M122").Value = -1845.5386
$ws.Range("N122").Value = -270094.33

$ws.Range("H130").Value = 54989.332
$ws.Range("J130").Value = 54989.332
$ws.Range("L130").Value = 54989.332
$ws.Range("N130").Value = -65029.332

$ws.Range("H136").Value = 56672.547
$ws.Range("J136").Value = 56672.547
$ws.Range("L136").Value = 56672.547
$ws.Range("N136").Value = -66872.54699999999

$ws.Range("H141").Value = 4429.6816
$ws.Range("I141").Value = 2535.6
$ws.Range("J141").Value = 8488.429
$ws.Range("K141").Value = 7606.799999999999
$ws.Range("L141").Value = 25465.287
$ws.Range("M141").Value = -2426.799999999999
$ws.Range("N141").Value = -35825.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21074.217
$ws.Range("I32").Value = 21626.059
$ws.Range("J32").Value = 9485.5
$ws.Range("K32").Value = 21626.059
$ws.Range("L32").Value = 9485.5
$ws.Range("M32").Value = -21339.059
$ws.Range("N32").Value = -10059.5

$ws.Range("H102").Value = 10971.45
$ws.Range("I102").Value = 994.94446
$ws.Range("J102").Value = 100760
$ws.Range("K102").Value = 994.94446
$ws.Range("L102").Value = 100760
$ws.Range("M102").Value = 627.05554
$ws.Range("N102").Value = -104004

$ws.Range("H110").Value = 1598.625
$ws.Range("I110").Value = 1508.2084
$ws.Range("J110").Value = 1869.875
$ws.Range("K110").Value = 1508.2084
$ws.Range("L110").Value = 1869.875
$ws.Range("M110").Value = 536.7916
$ws.Range("N110").Value = -5959.875

$ws.Range("H122").Value = 2509.1875
$ws.Range("I122").Value = 2745.9
$ws.Range("J122").Value = 2114.6667
$ws.Range("K122").Value = 8237.700000000001
$ws.Range("L122").Value = 6344.000100000001
$ws.Range("M122").Value = -5787.700000000001
$ws.Range("N122").Value = -11244.0001

$ws.Range("H129").Value = 33333
$ws.Range("J129").Value = 33333
$ws.Range("L129").Value = 33333
$ws.Range("N129").Value = -43333

$ws.Range("H131").Value = 51235.668
$ws.Range("J131").Value = 51235.668
$ws.Range("L131").Value = 51235.668
$ws.Range("N131").Value = -61315.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H94").Value = 695.23334
$ws.Range("I94").Value = 626.3913
$ws.Range("J94").Value = 921.4286
$ws.Range("K94").Value = 626.3913
$ws.Range("L94").Value = 921.4286
$ws.Range("M94").Value = -175.3913
$ws.Range("N94").Value = -1823.4286

$ws.Range("H95").Value = 42000
$ws.Range("J95").Value = 42000
$ws.Range("L95").Value = 42000
$ws.Range("N95").Value = -47492

$ws.Range("H130").Value = 47309.332
$ws.Range("J130").Value = 47309.332
$ws.Range("L130").Value = 47309.332
$ws.Range("N130").Value = -57349.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4700.0737
$ws.Range("I31").Value = 1758.4546
$ws.Range("J31").Value = 7473.6
$ws.Range("K31").Value = 1758.4546
$ws.Range("L31").Value = 7473.6
$ws.Range("M31").Value = -1463.4546
$ws.Range("N31").Value = -8063.6

$ws.Range("H34").Value = 4700.0737
$ws.Range("I34").Value = 1758.4546
$ws.Range("J34").Value = 7473.6
$ws.Range("K34").Value = 1758.4546
$ws.Range("L34").Value = 7473.6
$ws.Range("M34").Value = -1556.4546
$ws.Range("N34").Value = -7877.6

$ws.Range("H58").Value = 1843.4429
$ws.Range("I58").Value = 1620.8983
$ws.Range("J58").Value = 3037.0908
$ws.Range("K58").Value = 1620.8983
$ws.Range("L58").Value = 3037.0908
$ws.Range("M58").Value = -1417.8983
$ws.Range("N58").Value = -3443.0908

$ws.Range("H107").Value = 1132.6666
$ws.Range("J107").Value = 1045.75
$ws.Range("L107").Value = 1045.75
$ws.Range("N107").Value = -4885.75

$ws.Range("H109").Value = 34736.145
$ws.Range("J109").Value = 34736.145
$ws.Range("L109").Value = 34736.145
$ws.Range("N109").Value = -36816.145

$ws.Range("H136").Value = 1843.4429
$ws.Range("I136").Value = 1620.8983
$ws.Range("J136").Value = 3037.0908
$ws.Range("K136").Value = 4862.6949
$ws.Range("L136").Value = 9111.2724
$ws.Range("M136").Value = -2312.6949
$ws.Range("N136").Value = -14211.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2951.2195
$ws.Range("J11").Value = 2000
$ws.Range("L11").Value = 6000
$ws.Range("N11").Value = -6280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2775078.5
$ws.Range("J10").Value = 157
$ws.Range("L10").Value = 157
$ws.Range("N10").Value = -495

$ws.Range("H87").Value = 25714.285
$ws.Range("J87").Value = 25714.285
$ws.Range("L87").Value = 25714.285
$ws.Range("N87").Value = -28210.285

$ws.Range("H90").Value = 25714.285
$ws.Range("J90").Value = 25714.285
$ws.Range("L90").Value = 77142.855
$ws.Range("N90").Value = -89622.855

$ws.Range("H107").Value = 185212.19
$ws.Range("I107").Value = 333425.16
$ws.Range("J107").Value = 7356.6
$ws.Range("K107").Value = 333425.16
$ws.Range("L107").Value = 7356.6
$ws.Range("M107").Value = -331505.16
$ws.Range("N107").Value = -11196.6

$ws.Range("H122").Value = 2178.5715
$ws.Range("I122").Value = 2360
$ws.Range("J122").Value = 1725
$ws.Range("K122").Value = 7080
$ws.Range("L122").Value = 5175
$ws.Range("M122").Value = -4630
$ws.Range("N122").Value = -10075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8241.412
$ws.Range("J2").Value = 60001
$ws.Range("L2").Value = 60001
$ws.Range("N2").Value = -60225

$ws.Range("H16").Value = 1056.9642
$ws.Range("I16").Value = 1149.7916
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1149.7916
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -979.7916
$ws.Range("N16").Value = -840

$ws.Range("H46").Value = 4916.6665
$ws.Range("I46").Value = 4500
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 4500
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -4312
$ws.Range("N46").Value = -5376

$ws.Range("H93").Value = 2416.9167
$ws.Range("I93").Value = 2237.875
$ws.Range("J93").Value = 2775
$ws.Range("K93").Value = 2237.875
$ws.Range("L93").Value = 2775
$ws.Range("M93").Value = -989.875
$ws.Range("N93").Value = -5271

$ws.Range("H98").Value = 42494
$ws.Range("J98").Value = 42494
$ws.Range("L98").Value = 42494
$ws.Range("N98").Value = -48484

$ws.Range("H123").Value = 33284.668
$ws.Range("J123").Value = 33284.668
$ws.Range("L123").Value = 33284.668
$ws.Range("N123").Value = -43084.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 5002.3335
$ws.Range("J18").Value = 5002.3335
$ws.Range("L18").Value = 5002.3335
$ws.Range("N18").Value = -5348.3335

$ws.Range("H96").Value = 791
$ws.Range("I96").Value = 791
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 791
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 582
$ws.Range("N96").ClearContents()

$ws.Range("H109").Value = 32453
$ws.Range("J109").Value = 35943.6
$ws.Range("L109").Value = 35943.6
$ws.Range("N109").Value = -38717.6

$ws.Range("H127").Value = 41426.332
$ws.Range("J127").Value = 41426.332
$ws.Range("L127").Value = 41426.332
$ws.Range("N127").Value = -51346.332
